# Apply latest crypto price/volume snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.600.55"
$ws.Range("E2").Value = "  -1.49%  "

# Row 3
$ws.Range("D3").Value = "3.544.83"
$ws.Range("E3").Value = "  -2.87%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.08"
$ws.Range("E5").Value = "  -3.61%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "184.92"
$ws.Range("E6").Value = "  -4.58%  "

# Row 7
$ws.Range("D7").Value = "3.540.40"
$ws.Range("E7").Value = "  -2.83%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.615"
$ws.Range("E8").Value = "  -4.96%  "

# Row 9
$ws.Range("E9").Value = "  +0.04%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.180"
$ws.Range("E10").Value = "  -1.70%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.642"
$ws.Range("E11").Value = "  -4.64%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.29"
$ws.Range("E12").Value = "  -6.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000296"
$ws.Range("E13").Value = "  +1.04%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.42"
$ws.Range("E14").Value = "  -5.24%  "

# Row 15
$ws.Range("D15").Value = "4.122.45"
$ws.Range("E15").Value = "  -2.66%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.39"
$ws.Range("E16").Value = "  -3.37%  "

# Row 17
$ws.Range("D17").Value = "3.553.40"
$ws.Range("E17").Value = "  -2.68%  "

# Row 18
$ws.Range("D18").Value = "69.530.56"
$ws.Range("E18").Value = "  -1.64%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.39"
$ws.Range("E19").Value = "  -3.13%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.120"
$ws.Range("E20").Value = "  -1.02%  "

# Row 21
$ws.Range("E21").Value = "  -4.51%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "483.55"
$ws.Range("E22").Value = "  -1.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.25"
$ws.Range("E23").Value = "  +0.37%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.85"
$ws.Range("E24").Value = "  -8.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.33"
$ws.Range("E25").Value = "  -4.44%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "94.68"
$ws.Range("E26").Value = "  +3.55%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.33"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.91"
$ws.Range("E28").Value = "  -7.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.18"
$ws.Range("E29").Value = "  -4.44%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.27"
$ws.Range("E30").Value = "  -4.60%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.42"
$ws.Range("E31").Value = "  -4.57%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "66.31"
$ws.Range("E32").Value = "  +0.34%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.89"
$ws.Range("E33").Value = "  -3.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.113"
$ws.Range("E34").Value = "  -7.19%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "563.86"
$ws.Range("E35").Value = "  -9.95%  "

# Row 36
$ws.Range("E36").Value = "  +12.45%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "38.23"
$ws.Range("E37").Value = "  -5.14%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  -0.14%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.390"
$ws.Range("E39").Value = "  -5.37%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("E40").Value = "  -5.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.46"
$ws.Range("E41").Value = "  -3.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.13"
$ws.Range("E42").Value = "  +2.56%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  -10.08%  "

# Row 44
$ws.Range("B44").Value = "ThetaToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.96"
$ws.Range("E44").Value = "  -5.97%  "

# Row 45
$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "3.194.73"
$ws.Range("E45").Value = "  -3.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.43"
$ws.Range("E46").Value = "  +4.01%  "

# Row 47
$ws.Range("E47").Value = "  -5.34%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.38"
$ws.Range("E48").Value = "  +0.07%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.134"
$ws.Range("E49").Value = "  -3.46%  "

# Row 50
$ws.Range("E50").Value = "  +0.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.14"
$ws.Range("E51").Value = "  -4.54%  "
